$d = $word.ActiveDocument

# Replace the sponsorship deadline date "January 10" with "February 7th".
# The find text matches across the run boundary ("January " + "10") and the
# replacement collapses it into the single remaining run.
$d.Content.Find.Execute("January 10", $true, $false, $false, $false, $false,
                         $true, 1, $false, "February 7th", 2)
